$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New strings get appended to the shared-string table in the order they are
# first written, so "west" (row 37) is entered before "Millivolts" (row 33)
# to reproduce the source file's string order.
$ws.Range("E37").Value = "west"
$ws.Range("D33").Value = "Millivolts"

# New "Millivolts -> compass direction -> degrees" table in columns D:F,
# next to the existing Volts/compass table in A:B.
$ws.Range("D34").Value = 200
$ws.Range("E34").Value = "North"
$ws.Range("F34").Value = 0

$ws.Range("D35").Value = 800
$ws.Range("E35").Value = "East"
$ws.Range("F35").Value = 90

$ws.Range("D36").Value = 1400
$ws.Range("E36").Value = "South"
$ws.Range("F36").Value = 180

$ws.Range("D37").Value = 2000
$ws.Range("F37").Value = 270

$ws.Range("D38").Value = 2600
$ws.Range("E38").Value = "North"
$ws.Range("F38").Value = 360

# Match the saved selection / scroll position of the new table.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D33:E38").Select()
